# Decrement the "remaining days" counter (column E) for every data row.
# When the counter would drop to/below 0 (i.e. it was 1), the cycle resets:
# the counter goes back to 10 and the start date (column F) advances by 10 days.
# Row 36 has a malformed start-date value and is intentionally left untouched,
# matching the source data set which skips it as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
$skipRows = @(36)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) {
        continue
    }

    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value()
    if ($eVal -eq $null) {
        continue
    }

    if ($eVal -gt 1) {
        $eCell.Value = $eVal - 1
    } else {
        $fVal = $fCell.Value()
        $eCell.Value = 10
        $fCell.Value = $fVal + 10
    }
}
